$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$featuresText = "11 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, b#max-digit-skip-all-punctuation >= 7, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii"
$nnText = "Neural-Network"
$p1000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$p2000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$p3000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 3000"

$rows = @(
  @("20160426_092120", 0.924092409240924, 0.53, $p1000),
  @("20160426_093618", 0.917491749174917, 0.47, $p1000),
  @("20160426_095151", 0.920792079207921, 0.58, $p1000),
  @("20160426_100710", 0.920792079207921, 0.56, $p1000),
  @("20160426_102225", 0.920792079207921, 0.65, $p1000),
  @("20160426_110828", 0.920792079207921, 0.46, $p2000),
  @("20160426_113642", 0.920792079207921, 0.46, $p2000),
  @("20160426_120611", 0.914191419141914, 0.45, $p2000),
  @("20160426_123506", 0.920792079207921, 0.46, $p2000),
  @("20160426_130403", 0.917491749174917, 0.45, $p2000),
  @("20160426_134427", 0.917491749174917, 0.45, $p3000),
  @("20160426_142738", 0.920792079207921, 0.46, $p3000),
  @("20160426_151051", 0.914191419141914, 0.44, $p3000),
  @("20160426_155414", 0.917491749174917, 0.45, $p3000),
  @("20160426_163611", 0.920792079207921, 0.46, $p3000)
)

$r = 12
foreach ($row in $rows) {
    $time = $row[0]
    $classifyAcc = $row[1]
    $segmentAcc = $row[2]
    $paramText = $row[3]

    $ws.Cells.Item($r, 1).Value = $time
    $ws.Cells.Item($r, 2).Value = $featuresText
    $ws.Cells.Item($r, 3).Value = $featuresText
    $ws.Cells.Item($r, 4).Value = $featuresText
    $ws.Cells.Item($r, 5).Value = $nnText
    $ws.Cells.Item($r, 6).Value = $paramText
    $ws.Cells.Item($r, 7).Value = $nnText
    $ws.Cells.Item($r, 8).Value = $paramText
    $ws.Cells.Item($r, 9).Value = $nnText
    $ws.Cells.Item($r, 10).Value = $paramText
    $ws.Cells.Item($r, 11).Value = $classifyAcc
    $ws.Cells.Item($r, 12).Value = $segmentAcc

    $r = $r + 1
}
